# Adapt column header formatting to respective input file names (#7)
#
# - Renames the "*_old" header row (columns A:J) to "*_FV2410"
# - Renames the "*_new" header row (columns L:U) to "*_FV2504"
# - Wraps the used range A1:U58 in an Excel Table ("Table1")
# - Freezes the header row (row 1) in the active sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------

$oldToFv2410 = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
}

$newToFv2504 = @{
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $oldToFv2410.Keys) {
    $ws.Range($addr).Value = $oldToFv2410[$addr]
}

foreach ($addr in $newToFv2504.Keys) {
    $ws.Range($addr).Value = $newToFv2504[$addr]
}

# --- 2. Turn the used range into an Excel Table -----------------------------

$tableRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
